# The "Periodo Mora" (E16:E21) / "Valor Mora" (F16:F21) rows are
# refreshed: the previous account-statement periods are dropped and the
# new set is written in reverse (most-recent-first) order, per the
# commit "Elimina EC anteriores y se agregan nuevos, se modifica base
# de datos".
#
# Net effect on the data (row position / styling stays put, only the
# Periodo Mora + Valor Mora values change):
#   row 16: 1806 / 31249  ->  1902 / 28124
#   row 17: 1807 / 31249  ->  1810 / 31249
#   row 18: 1808 / 31249  ->  1809 / 31249
#   row 19: 1809 / 31249  ->  1808 / 31249
#   row 20: 1810 / 31249  ->  1807 / 31249
#   row 21: 1902 / 28124  ->  1806 / 31249

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1902", "1810", "1809", "1808", "1807", "1806")
$valores  = @(28124, 31249, 31249, 31249, 31249, 31249)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
